$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" for the second worker row (E17): 2506 -> 2508
$ws.Range("E17").Value = "2508"

# Update "Salario Basico" values (G16, G17): 877803 -> 1423500
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
